$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the slope comparison figures (start/end spawn temp) for row 7
$ws.Range("I7").Value = 5.99
$ws.Range("J7").Value = 4.7

# Update the active selection to I14
$ws.Range("I14").Select()
